# Roll the 90-day GSC export window forward by one day:
#   - drop 2025-10-28 (the oldest date)
#   - every other date shifts back one row (HTTPS URL counts follow it)
#   - append 2026-01-26 as the newest date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# New HTTPS-URL counts for rows 2..91 (column C), already shifted up by one
# row relative to the original data (row 91 keeps its prior value, matching
# the export - no new count was shown for the newest day).
$counts = @(
    90, 93, 92, 102, 113, 115, 107, 105, 100, 94,
    86, 83, 66, 54, 46, 43, 40, 37, 35, 30,
    29, 26, 25, 25, 26, 26, 25, 25, 27, 28,
    28, 27, 27, 27, 27, 27, 26, 25, 25, 25,
    26, 27, 27, 29, 29, 30, 30, 31, 31, 31,
    31, 31, 32, 32, 32, 32, 30, 31, 32, 30,
    28, 28, 28, 28, 29, 29, 28, 27, 27, 28,
    27, 27, 27, 27, 26, 26, 27, 26, 26, 25,
    25, 25, 25, 26, 25, 24, 23, 24, 24, 24
)

$startDate = Get-Date -Year 2025 -Month 10 -Day 29

for ($i = 0; $i -lt 90; $i++) {
    $r = $i + 2

    $dateText = $startDate.AddDays($i).ToString("yyyy-MM-dd")

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dateText
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 3).Value = $counts[$i]
}
